$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(175).Insert()

$ws.Cells.Item(175, 1).Value = 3
$ws.Cells.Item(175, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(175, 3).Value = "Coquimbo"
$ws.Cells.Item(175, 4).Value = 45029
$ws.Cells.Item(175, 5).Value = 5
$ws.Cells.Item(175, 6).Value = 100112052
$ws.Cells.Item(175, 7).Value = "Albahaca"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 130
$ws.Cells.Item(175, 11).Value = 4000
$ws.Cells.Item(175, 12).Value = 4500
$ws.Cells.Item(175, 13).Value = 4269
$ws.Cells.Item(175, 14).Value = "$/docena de matas"
$ws.Cells.Item(175, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(175, 16).Value = 712
$ws.Cells.Item(175, 17).Value = 6
$ws.Cells.Item(175, 18).Value = "Hortaliza"
